$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: UID=2, same Source/Destination plate info, Destination Well = A2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("F3").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = "A2"
$ws.Range("H3").Value = 2875
$ws.Range("I3").Value = $ws.Range("I2").Value2

# Row 4: UID=3, same Source/Destination plate info, Destination Well = A3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = $ws.Range("B2").Value2
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("E4").Value = $ws.Range("E2").Value2
$ws.Range("F4").Value = $ws.Range("F2").Value2
$ws.Range("G4").Value = "A3"
$ws.Range("H4").Value = 2875
$ws.Range("I4").Value = $ws.Range("I2").Value2
